$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed "Price" (D) / "Volume(1h)" (E) snapshot for the crypto list.
# Every value in these columns is stored as plain text in the sheet
# (<c t="inlineStr">), so cells whose new reading happens to look like a plain
# number (IsNumeric = $true below) are pre-formatted as text first -- otherwise
# Excel would silently coerce e.g. "540.37" into a numeric value.
$updates = @(
    @{ Cell = "D2"; Value = "63.730.06"; IsNumeric = $false }
    @{ Cell = "E2"; Value = "  +0.51%  "; IsNumeric = $false }
    @{ Cell = "D3"; Value = "3.090.15"; IsNumeric = $false }
    @{ Cell = "E3"; Value = "  -1.04%  "; IsNumeric = $false }
    @{ Cell = "E4"; Value = "  +0.07%  "; IsNumeric = $false }
    @{ Cell = "D5"; Value = "540.37"; IsNumeric = $true }
    @{ Cell = "E5"; Value = "  -3.11%  "; IsNumeric = $false }
    @{ Cell = "D6"; Value = "136.67"; IsNumeric = $true }
    @{ Cell = "E6"; Value = "  -2.23%  "; IsNumeric = $false }
    @{ Cell = "D7"; Value = "1.00"; IsNumeric = $true }
    @{ Cell = "E7"; Value = "  +0.03%  "; IsNumeric = $false }
    @{ Cell = "D8"; Value = "3.083.12"; IsNumeric = $false }
    @{ Cell = "E8"; Value = "  -1.02%  "; IsNumeric = $false }
    @{ Cell = "E9"; Value = "  -0.36%  "; IsNumeric = $false }
    @{ Cell = "E10"; Value = "  -2.83%  "; IsNumeric = $false }
    @{ Cell = "D11"; Value = "6.28"; IsNumeric = $true }
    @{ Cell = "E11"; Value = "  -6.01%  "; IsNumeric = $false }
    @{ Cell = "D12"; Value = "0.458"; IsNumeric = $true }
    @{ Cell = "E12"; Value = "  -0.61%  "; IsNumeric = $false }
    @{ Cell = "E13"; Value = "  +3.69%  "; IsNumeric = $false }
    @{ Cell = "D14"; Value = "34.79"; IsNumeric = $true }
    @{ Cell = "E14"; Value = "  -2.04%  "; IsNumeric = $false }
    @{ Cell = "D15"; Value = "3.589.26"; IsNumeric = $false }
    @{ Cell = "E15"; Value = "  -1.11%  "; IsNumeric = $false }
    @{ Cell = "D16"; Value = "63.758.83"; IsNumeric = $false }
    @{ Cell = "E16"; Value = "  +0.51%  "; IsNumeric = $false }
    @{ Cell = "E17"; Value = "  +0.16%  "; IsNumeric = $false }
    @{ Cell = "D18"; Value = "3.087.85"; IsNumeric = $false }
    @{ Cell = "E18"; Value = "  -0.95%  "; IsNumeric = $false }
    @{ Cell = "D19"; Value = "6.70"; IsNumeric = $true }
    @{ Cell = "E19"; Value = "  -0.88%  "; IsNumeric = $false }
    @{ Cell = "D20"; Value = "488.86"; IsNumeric = $true }
    @{ Cell = "E20"; Value = "  -4.34%  "; IsNumeric = $false }
    @{ Cell = "D21"; Value = "13.49"; IsNumeric = $true }
    @{ Cell = "E21"; Value = "  -1.36%  "; IsNumeric = $false }
    @{ Cell = "D22"; Value = "0.702"; IsNumeric = $true }
    @{ Cell = "E22"; Value = "  -1.59%  "; IsNumeric = $false }
    @{ Cell = "D23"; Value = "7.19"; IsNumeric = $true }
    @{ Cell = "E23"; Value = "  -1.94%  "; IsNumeric = $false }
    @{ Cell = "D24"; Value = "79.86"; IsNumeric = $true }
    @{ Cell = "E24"; Value = "  +1.85%  "; IsNumeric = $false }
    @{ Cell = "D25"; Value = "12.27"; IsNumeric = $true }
    @{ Cell = "E25"; Value = "  -1.78%  "; IsNumeric = $false }
    @{ Cell = "E26"; Value = "  +0.03%  "; IsNumeric = $false }
    @{ Cell = "E27"; Value = "  -1.63%  "; IsNumeric = $false }
    @{ Cell = "D28"; Value = "8.32"; IsNumeric = $true }
    @{ Cell = "E28"; Value = "  -0.46%  "; IsNumeric = $false }
    @{ Cell = "D29"; Value = "0.998"; IsNumeric = $true }
    @{ Cell = "E29"; Value = "  -0.06%  "; IsNumeric = $false }
    @{ Cell = "D30"; Value = "26.27"; IsNumeric = $true }
    @{ Cell = "E30"; Value = "  -1.04%  "; IsNumeric = $false }
    @{ Cell = "E31"; Value = "  -3.31%  "; IsNumeric = $false }
    @{ Cell = "E32"; Value = "  -0.49%  "; IsNumeric = $false }
    @{ Cell = "E33"; Value = "  -5.29%  "; IsNumeric = $false }
    @{ Cell = "D34"; Value = "56.98"; IsNumeric = $true }
    @{ Cell = "E34"; Value = "  -4.79%  "; IsNumeric = $false }
    @{ Cell = "D35"; Value = "5.52"; IsNumeric = $true }
    @{ Cell = "E35"; Value = "  +4.76%  "; IsNumeric = $false }
    @{ Cell = "D36"; Value = "495.67"; IsNumeric = $true }
    @{ Cell = "E36"; Value = "  -7.58%  "; IsNumeric = $false }
    @{ Cell = "D37"; Value = "6.07"; IsNumeric = $true }
    @{ Cell = "E37"; Value = "  +0.97%  "; IsNumeric = $false }
    @{ Cell = "D38"; Value = "3.294.50"; IsNumeric = $false }
    @{ Cell = "E38"; Value = "  +7.01%  "; IsNumeric = $false }
    @{ Cell = "E39"; Value = "  -4.17%  "; IsNumeric = $false }
    @{ Cell = "D40"; Value = "0.0801"; IsNumeric = $true }
    @{ Cell = "E40"; Value = "  -0.14%  "; IsNumeric = $false }
    @{ Cell = "E41"; Value = "  -2.86%  "; IsNumeric = $false }
    @{ Cell = "D42"; Value = "8.16"; IsNumeric = $true }
    @{ Cell = "E42"; Value = "  -0.15%  "; IsNumeric = $false }
    @{ Cell = "D43"; Value = "2.65"; IsNumeric = $true }
    @{ Cell = "E43"; Value = "  -3.69%  "; IsNumeric = $false }
    @{ Cell = "E44"; Value = "  +0.26%  "; IsNumeric = $false }
    @{ Cell = "D46"; Value = "2.08"; IsNumeric = $true }
    @{ Cell = "E46"; Value = "  -0.67%  "; IsNumeric = $false }
    @{ Cell = "D47"; Value = "0.0₃0538"; IsNumeric = $false }
    @{ Cell = "E47"; Value = "  +4.24%  "; IsNumeric = $false }
    @{ Cell = "D48"; Value = "24.98"; IsNumeric = $true }
    @{ Cell = "E48"; Value = "  +2.01%  "; IsNumeric = $false }
    @{ Cell = "D49"; Value = "121.73"; IsNumeric = $true }
    @{ Cell = "E49"; Value = "  -0.89%  "; IsNumeric = $false }
    @{ Cell = "E50"; Value = "  +1.53%  "; IsNumeric = $false }
    @{ Cell = "E51"; Value = "  -3.65%  "; IsNumeric = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.IsNumeric) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
    if ($u.IsNumeric) {
        $cell.Style = "Normal"
    }
}
